# Generate Report for Handback
# Refresh the handoff/handback timestamps and priority/status values that
# get written each time the handback report is (re)generated.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview.Range("G2").Value = "2016-08-30 22:17:07"
$wsOverview.Range("G5").Value = "2016-08-30 22:17:07"

# --- zh-cn sheet ---
# Priority (E)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"

# Correspond Handoff Datetime (H)
$wsZhCn.Range("H2").Value = "2016-08-30 22:16:57"
$wsZhCn.Range("H5").Value = "2016-08-30 22:16:57"

# Correspond Handback DateTime (K)
$wsZhCn.Range("K2").Value = "2016-08-30 22:17:28"
$wsZhCn.Range("K5").Value = "2016-08-30 22:17:28"

# --- de-de sheet ---
# Priority (E)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"

# Correspond Handoff Datetime (H) -- shares the same value as Overview!G
$wsDeDe.Range("H2").Value = "2016-08-30 22:17:07"
$wsDeDe.Range("H5").Value = "2016-08-30 22:17:07"

# Correspond Handback DateTime (K)
$wsDeDe.Range("K2").Value = "2016-08-30 22:17:37"
$wsDeDe.Range("K5").Value = "2016-08-30 22:17:37"
